# Update "想去人数" (wanted-to-go headcount) values in column F across sheets
# for the gh-pages data refresh generated at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 506
$ws.Range("F4").Value = 486
$ws.Range("F5").Value = 902
$ws.Range("F7").Value = 866
$ws.Range("F8").Value = 679
$ws.Range("F9").Value = 136
$ws.Range("F11").Value = 67
$ws.Range("F13").Value = 231
$ws.Range("F14").Value = 526
$ws.Range("F16").Value = 1256
$ws.Range("F17").Value = 107
$ws.Range("F18").Value = 988
$ws.Range("F19").Value = 2732
$ws.Range("F20").Value = 1205
$ws.Range("F21").Value = 623
$ws.Range("F22").Value = 151
$ws.Range("F23").Value = 1215
$ws.Range("F25").Value = 928
$ws.Range("F26").Value = 101
$ws.Range("F27").Value = 1251

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 496
$ws.Range("F8").Value = 34

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 709

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 709
$ws.Range("F3").Value = 506
$ws.Range("F6").Value = 486
$ws.Range("F7").Value = 496
$ws.Range("F8").Value = 496
$ws.Range("F12").Value = 902
$ws.Range("F14").Value = 866
$ws.Range("F15").Value = 679
$ws.Range("F16").Value = 136
$ws.Range("F19").Value = 34
$ws.Range("F22").Value = 67
$ws.Range("F25").Value = 231
$ws.Range("F26").Value = 526
$ws.Range("F28").Value = 1256
$ws.Range("F29").Value = 107
$ws.Range("F30").Value = 988
$ws.Range("F31").Value = 2732
$ws.Range("F32").Value = 1205
$ws.Range("F33").Value = 623
$ws.Range("F34").Value = 152
$ws.Range("F35").Value = 1215
$ws.Range("F38").Value = 928
$ws.Range("F39").Value = 101
$ws.Range("F40").Value = 1251

